$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 header edits ----
$ws.Range("K1").Value = "Precio/m²Características"

# Q1/R1 are new header cells - copy the header style (bold/border/centered)
# from an existing header cell (P1) before writing their text.
$ws.Range("P1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "Características Básicas"
$ws.Range("P1").Copy($ws.Range("R1"))
$ws.Range("R1").Value = "Más Características"

# ---- Row 2 (new data row) ----
$ws.Range("A2").Value = 105773947
$ws.Range("B2").Value = "Alquiler"
$ws.Range("C2").Value = "Alquiler de Ático en Peñas Negras, 12"
$ws.Range("D2").Value = "(`'`nPeñas Negras, 12`n`',)"
$ws.Range("E2").Value = "(`'`nBarrio Torreagüera`n`',)"
$ws.Range("F2").Value = "(`'`nDistrito Pedanías Este`n`',)"
$ws.Range("G2").Value = "(`'`nMurcia`n`',)"
$ws.Range("H2").Value = "`nÁrea de Murcia, Murcia`n"
$ws.Range("I2").Value = "700€/mes"
$ws.Range("L2").Value = 105773947
$ws.Range("M2").Value = "Profesional"
$ws.Range("N2").Value = "CITYSOL MURCIA"
$ws.Range("O2").Value = "Anuncio actualizado el 20 de agosto"
$ws.Range("P2").Value = "[]"
$ws.Range("Q2").Value = "['128 m² construidos', '3 habitaciones', '2 baños', 'Terraza y balcón', 'Plaza de garaje incluida en el precio', 'Segunda mano/buen estado', 'Armarios empotrados', 'Trastero', 'Orientación este, oeste', 'Construido en 2010', 'Cocina equipada y casa sin amueblar', 'Planta 3ª exterior', 'Con ascensor']"
$ws.Range("R2").Value = "[]"

# ---- Row 3 (new data row) ----
# A3 and L3 are digit-strings that must stay TEXT (not be coerced to numbers).
$ws.Range("A3").Value = "'37324696"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Alquiler"
$ws.Range("C3").Value = "Alquiler de Piso en calle Jacobo de las Leyes"
$ws.Range("D3").Value = "(`'`nCalle Jacobo de las Leyes`n`',)"
$ws.Range("E3").Value = "(`'`nBarrio La Fama`n`',)"
$ws.Range("F3").Value = "(`'`nDistrito Centro`n`',)"
$ws.Range("G3").Value = "(`'`nMurcia`n`',)"
$ws.Range("H3").Value = "`nÁrea de Murcia, Murcia`n"
$ws.Range("I3").Value = "750€/mes"
$ws.Range("J3").Value = "Fianza de 1 mes"
$ws.Range("L3").Value = "'2530"
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Value = "Profesional"
$ws.Range("N3").Value = "The Simple Rent"
$ws.Range("O3").Value = "Anuncio actualizado el 20 de agosto"
$ws.Range("P3").Value = "N/A"
$ws.Range("Q3").Value = "['90 m² construidos, 85 m² útiles', '3 habitaciones', '1 baño', 'Balcón', 'Segunda mano/buen estado', 'Armarios empotrados', 'Trastero', 'Amueblado y cocina equipada', 'Planta 3ª exterior', 'Sin ascensor']"
$ws.Range("R3").Value = "[]"
